$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 18.714285
$ws.Range("I6").Value2 = 4.2
$ws.Range("K6").Value2 = 12.6
$ws.Range("M6").Value2 = 99.40000000000001
$ws.Range("H32").Value2 = 9221.444
$ws.Range("I32").Value2 = 8248.25
$ws.Range("K32").Value2 = 8248.25
$ws.Range("M32").Value2 = -7922.25
$ws.Range("H33").Value2 = 166.8
$ws.Range("I33").Value2 = 166.8
$ws.Range("K33").Value2 = 166.8
$ws.Range("M33").Value2 = 62.19999999999999
$ws.Range("H42").Value2 = 58.333332
$ws.Range("I42").Value2 = 50
$ws.Range("J42").Value2 = 66.666664
$ws.Range("K42").Value2 = 150
$ws.Range("L42").Value2 = 199.999992
$ws.Range("M42").Value2 = 80
$ws.Range("N42").Value2 = -659.999992
$ws.Range("H81").Value2 = 0
$ws.Range("I81").Value2 = 0
$ws.Range("K81").Value2 = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value2 = 0
$ws.Range("I84").Value2 = 0
$ws.Range("K84").Value2 = 0
$ws.Range("M84").ClearContents()
$ws.Range("H121").Value2 = 0
$ws.Range("J121").Value2 = 0
$ws.Range("L121").Value2 = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value2 = 6148.5
$ws.Range("I132").Value2 = 2659.6667
$ws.Range("K132").Value2 = 7979.000100000001
$ws.Range("M132").Value2 = -5449.000100000001
$ws.Range("H135").Value2 = 2987
$ws.Range("I135").Value2 = 2987
$ws.Range("K135").Value2 = 26883
$ws.Range("M135").Value2 = -24348
$ws.Range("H137").Value2 = 1744
$ws.Range("I137").Value2 = 1590.4
$ws.Range("J137").Value2 = 2000
$ws.Range("K137").Value2 = 4771.200000000001
$ws.Range("L137").Value2 = 6000
$ws.Range("M137").Value2 = -2221.200000000001
$ws.Range("N137").Value2 = -11100
$ws.Range("H138").Value2 = 3042.2
$ws.Range("J138").Value2 = 3999
$ws.Range("L138").Value2 = 11997
$ws.Range("N138").Value2 = -22277

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value2 = 2475
$ws.Range("I63").Value2 = 2475
$ws.Range("K63").Value2 = 2475
$ws.Range("M63").Value2 = -1789
$ws.Range("H66").Value2 = 2475
$ws.Range("I66").Value2 = 2475
$ws.Range("K66").Value2 = 12375
$ws.Range("M66").Value2 = -8943

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 4655.8
$ws.Range("I86").Value2 = 2511.6
$ws.Range("K86").Value2 = 2511.6
$ws.Range("M86").Value2 = -1388.6
$ws.Range("H89").Value2 = 4655.8
$ws.Range("I89").Value2 = 2511.6
$ws.Range("K89").Value2 = 12558
$ws.Range("M89").Value2 = -6942
$ws.Range("H94").Value2 = 4037
$ws.Range("I94").Value2 = 3555.5
$ws.Range("K94").Value2 = 3555.5
$ws.Range("M94").Value2 = -3104.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1001
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 1001
$ws.Range("K22").Value2 = 0
$ws.Range("L22").Value2 = 1001
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value2 = -1701
$ws.Range("H82").Value2 = 0
$ws.Range("I82").Value2 = 0
$ws.Range("K82").Value2 = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("M85").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value2 = 3000
$ws.Range("J3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("N3").ClearContents()
$ws.Range("H14").Value2 = 921.2
$ws.Range("I14").Value2 = 921.2
$ws.Range("K14").Value2 = 2763.6
$ws.Range("M14").Value2 = -2590.6
$ws.Range("H44").Value2 = 0
$ws.Range("I44").Value2 = 0
$ws.Range("K44").Value2 = 0
$ws.Range("M44").ClearContents()
$ws.Range("H52").Value2 = 881.5
$ws.Range("J52").Value2 = 500
$ws.Range("L52").Value2 = 1500
$ws.Range("N52").Value2 = -2032
$ws.Range("H115").Value2 = 31
$ws.Range("J115").Value2 = 31
$ws.Range("L115").Value2 = 93
$ws.Range("N115").Value2 = -2443
$ws.Range("H122").Value2 = 3633.8
$ws.Range("J122").Value2 = 3633.8
$ws.Range("L122").Value2 = 32704.2
$ws.Range("N122").Value2 = -37604.2
$ws.Range("H131").Value2 = 2731.6667
$ws.Range("I131").Value2 = 2600
$ws.Range("J131").Value2 = 2995
$ws.Range("K131").Value2 = 7800
$ws.Range("L131").Value2 = 8985
$ws.Range("M131").Value2 = -2760
$ws.Range("N131").Value2 = -19065
$ws.Range("H138").Value2 = 0
$ws.Range("I138").Value2 = 0
$ws.Range("K138").Value2 = 0
$ws.Range("M138").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2023.7778
$ws.Range("I80").Value2 = 1402.3334
$ws.Range("K80").Value2 = 1402.3334
$ws.Range("M80").Value2 = -404.3334
$ws.Range("H83").Value2 = 2023.7778
$ws.Range("I83").Value2 = 1402.3334
$ws.Range("K83").Value2 = 7011.666999999999
$ws.Range("M83").Value2 = -2019.666999999999
$ws.Range("H107").Value2 = 1199
$ws.Range("I107").Value2 = 0
$ws.Range("J107").Value2 = 1199
$ws.Range("K107").Value2 = 0
$ws.Range("L107").Value2 = 1199
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value2 = -5039
$ws.Range("H109").Value2 = 45142.5
$ws.Range("J109").Value2 = 45142.5
$ws.Range("L109").Value2 = 45142.5
$ws.Range("N109").Value2 = -47222.5
$ws.Range("H132").Value2 = 2502.5454
$ws.Range("I132").Value2 = 2502.5454
$ws.Range("K132").Value2 = 7507.6362
$ws.Range("M132").Value2 = -4977.6362

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1584.25
$ws.Range("I16").Value2 = 1667.7142
$ws.Range("K16").Value2 = 1667.7142
$ws.Range("M16").Value2 = -1497.7142
$ws.Range("H46").Value2 = 748.1429000000001
$ws.Range("J46").Value2 = 696
$ws.Range("L46").Value2 = 696
$ws.Range("N46").Value2 = -1072
$ws.Range("H62").Value2 = 21633.334
$ws.Range("I62").Value2 = 20000
$ws.Range("J62").Value2 = 24900
$ws.Range("K62").Value2 = 20000
$ws.Range("L62").Value2 = 24900
$ws.Range("M62").Value2 = -19376
$ws.Range("N62").Value2 = -26148
$ws.Range("H65").Value2 = 21633.334
$ws.Range("I65").Value2 = 20000
$ws.Range("J65").Value2 = 24900
$ws.Range("K65").Value2 = 60000
$ws.Range("L65").Value2 = 74700
$ws.Range("M65").Value2 = -56880
$ws.Range("N65").Value2 = -80940

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("N66").ClearContents()
$ws.Range("H126").Value2 = 4171.7144
$ws.Range("I126").Value2 = 4125.75
$ws.Range("J126").Value2 = 4233
$ws.Range("K126").Value2 = 12377.25
$ws.Range("L126").Value2 = 12699
$ws.Range("M126").Value2 = -9907.25
$ws.Range("N126").Value2 = -17639
